$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting (avoid Excel auto-converting numeric-looking strings)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "56.993.52"
$ws.Range("E2").Value = "  +8.43%  "
$ws.Range("D3").Value = "3.246.52"
$ws.Range("E3").Value = "  +4.20%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "394.71"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").Value = "107.86"
$ws.Range("E6").Value = "  +3.89%  "
$ws.Range("D7").Value = "3.244.73"
$ws.Range("E7").Value = "  +4.24%  "
$ws.Range("D8").Value = "0.565"
$ws.Range("E8").Value = "  +4.66%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "0.619"
$ws.Range("E10").Value = "  +2.69%  "
$ws.Range("D11").Value = "39.13"
$ws.Range("E11").Value = "  +2.49%  "
$ws.Range("D12").Value = "0.0977"
$ws.Range("E12").Value = "  +13.32%  "
$ws.Range("E13").Value = "  +1.96%  "
$ws.Range("D14").Value = "3.768.52"
$ws.Range("E14").Value = "  +4.55%  "
$ws.Range("D15").Value = "8.12"
$ws.Range("E15").Value = "  +3.54%  "
$ws.Range("D16").Value = "19.01"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("D17").Value = "3.239.68"
$ws.Range("E17").Value = "  +3.92%  "
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("D19").Value = "10.65"
$ws.Range("E19").Value = "  -2.10%  "
$ws.Range("D20").Value = "56.923.52"
$ws.Range("E20").Value = "  +8.51%  "
$ws.Range("E21").Value = "  +2.63%  "
$ws.Range("E22").Value = "  +9.01%  "
$ws.Range("D23").Value = "13.04"
$ws.Range("E23").Value = "  +2.57%  "
$ws.Range("D24").Value = "299.13"
$ws.Range("E24").Value = "  +11.42%  "
$ws.Range("D25").Value = "73.92"
$ws.Range("E25").Value = "  +4.22%  "
$ws.Range("E26").Value = "  -2.50%  "
$ws.Range("D27").Value = "27.93"
$ws.Range("E27").Value = "  +1.36%  "
$ws.Range("D28").Value = "7.90"
$ws.Range("E28").Value = "  -1.95%  "
$ws.Range("D29").Value = "4.38"
$ws.Range("D30").Value = "0.169"
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("D31").Value = "7.24"
$ws.Range("E31").Value = "  -2.28%  "
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("D34").Value = "10.97"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").Value = "37.85"
$ws.Range("E35").Value = "  +2.88%  "
$ws.Range("D36").Value = "0.0484"
$ws.Range("E36").Value = "  -1.71%  "
$ws.Range("E37").Value = "  +1.65%  "
$ws.Range("D38").Value = "51.72"
$ws.Range("E38").Value = "  +3.52%  "
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "3.07"
$ws.Range("E40").Value = "  +14.20%  "
$ws.Range("D41").Value = "3.47"
$ws.Range("E41").Value = "  +1.82%  "
$ws.Range("D42").Value = "133.83"
$ws.Range("E42").Value = "  +2.80%  "
$ws.Range("E43").Value = "  +2.10%  "

# Row 44/45 swapped: Celestia <-> Stellar
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "0.120"
$ws.Range("E44").Value = "  +3.13%  "
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").Value = "17.13"
$ws.Range("E45").Value = "  +0.82%  "

$ws.Range("E46").Value = "  -3.31%  "
$ws.Range("E47").Value = "  -2.46%  "
$ws.Range("D48").Value = "21.94"
$ws.Range("E48").Value = "  -0.89%  "
$ws.Range("D49").Value = "2.140.48"
$ws.Range("E49").Value = "  +2.80%  "
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("E51").Value = "  -2.69%  "
